$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cxcl13"
$ws.Range("C2").Value = "Cxcr3"
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.767740666666666
$ws.Range("H2").Value = 23.303222
$ws.Range("I2").Value = 0.9880684546028148
$ws.Range("J2").Value = 0.9880684546028148
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.7166990000000001
$ws.Range("N2").Value = 2.150097
$ws.Range("O2").Value = 0.2276207788704612
$ws.Range("P2").Value = 0.2276207788704611
$ws.Range("Q2").Value = 5.567131968059333
$ws.Range("R2").Value = 50.104187712534
$ws.Range("S2").Value = 0.2249049112140256
$ws.Range("T2").Value = 0.2249049112140256

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cxcl13"
$ws.Range("C3").Value = "Cxcr3"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.767740666666666
$ws.Range("H3").Value = 23.303222
$ws.Range("I3").Value = 0.9880684546028148
$ws.Range("J3").Value = 0.9880684546028148
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.431954666666666
$ws.Range("N3").Value = 7.295864
$ws.Range("O3").Value = 0.7723792211295388
$ws.Range("P3").Value = 0.7723792211295388
$ws.Range("Q3").Value = 18.89079316375644
$ws.Range("R3").Value = 170.017138473808
$ws.Range("S3").Value = 0.7631635433887892
$ws.Range("T3").Value = 0.7631635433887892

$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Cxcl13"
$ws.Range("C4").Value = "Cxcr3"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.09380033333333333
$ws.Range("H4").Value = 0.281401
$ws.Range("I4").Value = 0.01193154539718528
$ws.Range("J4").Value = 0.01193154539718528
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.7166990000000001
$ws.Range("N4").Value = 2.150097
$ws.Range("O4").Value = 0.2276207788704612
$ws.Range("P4").Value = 0.2276207788704611
$ws.Range("Q4").Value = 0.06722660509966667
$ws.Range("R4").Value = 0.6050394458970001
$ws.Range("S4").Value = 0.002715867656435579
$ws.Range("T4").Value = 0.002715867656435578

$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Cxcl13"
$ws.Range("C5").Value = "Cxcr3"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.09380033333333333
$ws.Range("H5").Value = 0.281401
$ws.Range("I5").Value = 0.01193154539718528
$ws.Range("J5").Value = 0.01193154539718528
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.431954666666666
$ws.Range("N5").Value = 7.295864
$ws.Range("O5").Value = 0.7723792211295388
$ws.Range("P5").Value = 0.7723792211295388
$ws.Range("Q5").Value = 0.2281181583848889
$ws.Range("R5").Value = 2.053063425464
$ws.Range("S5").Value = 0.0092156777407497
$ws.Range("T5").Value = 0.0092156777407497
